# Check whether everything works with a new field
#
# Rows 6-9 currently show the "Edit by JS" long comment
# ("My long text\nsome new line.") in column B and the "ReadOnly" approval
# note ("Yes, it's approved.") in column D. Flip them around: column B
# should show the approval note and column D should show the long
# comment. Row 10's column D also picks up the long-comment
# value/formatting (row 10's column B already has it and is untouched).
#
# There's no native "swap"/"move" on a Range, so two fixed reference
# cells (one already holding each of the two target value+format
# combinations) are captured up front into off-sheet staging cells
# first, and every destination cell is then overwritten from the
# appropriate staged copy. Capturing both source patterns before any
# destination is touched avoids clobbering a source cell before it has
# been read (row 7's column D starts out blank, so it can't be used as
# its own source).
#
# xlPasteValues (-4163) copies the literal value across (re-using the
# existing shared-string entry rather than creating a duplicate) and
# xlPasteFormats (-4122) copies the number format/font/alignment, so
# doing both in sequence reproduces the source cell's full value+style
# on the destination.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues  = -4163

function Copy-CellFull {
    param($srcAddr, $dstAddr)
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteValues) | Out-Null
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# Stage the two value+format combinations from known-good source cells
# before mutating anything: "commentStage" = the long multi-line
# comment (currently column B's style/value), "approvedStage" = the
# approval note (currently column D's style/value on rows 6, 8 and 9).
$commentStage  = "Y1"
$approvedStage = "Z1"
Copy-CellFull "B6" $commentStage
Copy-CellFull "D6" $approvedStage

# Column B, rows 6-9: now shows the approval note.
Copy-CellFull $approvedStage "B6"
Copy-CellFull $approvedStage "B7"
Copy-CellFull $approvedStage "B8"
Copy-CellFull $approvedStage "B9"

# Column D, rows 6-10: now shows the long multi-line comment.
Copy-CellFull $commentStage "D6"
Copy-CellFull $commentStage "D7"
Copy-CellFull $commentStage "D8"
Copy-CellFull $commentStage "D9"
Copy-CellFull $commentStage "D10"

# Tidy up the staging cells.
$ws.Range("Y1:Z1").Clear() | Out-Null

# Leave the selection on the column D range that was just edited.
$ws.Range("D6:D10").Select() | Out-Null
